# Update "想去人数" (F) and "最低票价" (G) columns on both the "展览" and
# "全部类型" worksheets, which carry duplicated data tables.
#
# Map of row -> (new F value, new G value or $null to leave unchanged)
$updates = @{
    2  = @{ F = 3427; G = $null }
    3  = @{ F = 153;  G = 65 }
    4  = @{ F = 229;  G = 58 }
    5  = @{ F = 1764; G = 70 }
    6  = @{ F = 1669; G = 65 }
    7  = @{ F = 482;  G = 55 }
    13 = @{ F = 238;  G = $null }
    15 = @{ F = 56;   G = $null }
    16 = @{ F = 245;  G = $null }
    17 = @{ F = 236;  G = $null }
    20 = @{ F = 22;   G = $null }
    21 = @{ F = 66;   G = $null }
    22 = @{ F = 125;  G = $null }
    25 = @{ F = 296;  G = $null }
    30 = @{ F = 587;  G = $null }
    31 = @{ F = 2395; G = $null }
    35 = @{ F = 589;  G = $null }
    39 = @{ F = 364;  G = $null }
    41 = @{ F = 557;  G = $null }
}

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $vals = $updates[$row]
        $ws.Range("F$row").Value = $vals.F
        if ($null -ne $vals.G) {
            $ws.Range("G$row").Value = $vals.G
        }
    }
}
